# Auto update Excel log: append new sensor readings to PIR, Humidity, and Temperature sheets
$wb = $excel.ActiveWorkbook

function Set-TextRow($sheet, $row, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $sheet.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

$pir = $wb.Worksheets.Item("PIR")
$pirRange = $pir.Range("A204:F223")
$pirRange.NumberFormat = "@"
Set-TextRow $pir 204 @("2026-01-28", "12:19:20", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 205 @("2026-01-28", "12:19:20", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 206 @("2026-01-28", "12:19:22", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 207 @("2026-01-28", "12:19:24", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 208 @("2026-01-28", "12:19:25", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 209 @("2026-01-28", "12:19:27", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 210 @("2026-01-28", "12:19:30", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 211 @("2026-01-28", "12:19:32", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 212 @("2026-01-28", "12:19:34", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 213 @("2026-01-28", "12:19:35", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 214 @("2026-01-28", "12:19:37", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 215 @("2026-01-28", "12:19:39", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 216 @("2026-01-28", "12:19:42", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 217 @("2026-01-28", "12:19:47", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 218 @("2026-01-28", "12:19:53", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 219 @("2026-01-28", "12:19:57", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 220 @("2026-01-28", "12:20:02", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 221 @("2026-01-28", "12:20:07", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 222 @("2026-01-28", "12:20:12", "12:00", "Bathroom", "No Motion", "Inactive")
Set-TextRow $pir 223 @("2026-01-28", "12:20:18", "12:00", "Bathroom", "No Motion", "Inactive")
$pirRange.Style = "Normal"

$humidity = $wb.Worksheets.Item("Humidity")
$humidityRange = $humidity.Range("A191:F210")
$humidityRange.NumberFormat = "@"
Set-TextRow $humidity 191 @("2026-01-28", "12:19:19", "12:00", "Bathroom", "86.8%", "Active")
Set-TextRow $humidity 192 @("2026-01-28", "12:19:21", "12:00", "Bathroom", "86.8%", "Active")
Set-TextRow $humidity 193 @("2026-01-28", "12:19:23", "12:00", "Bathroom", "87.7%", "Active")
Set-TextRow $humidity 194 @("2026-01-28", "12:19:24", "12:00", "Bathroom", "86.8%", "Active")
Set-TextRow $humidity 195 @("2026-01-28", "12:19:26", "12:00", "Bathroom", "87.7%", "Active")
Set-TextRow $humidity 196 @("2026-01-28", "12:19:28", "12:00", "Bathroom", "86.9%", "Active")
Set-TextRow $humidity 197 @("2026-01-28", "12:19:29", "12:00", "Bathroom", "87.7%", "Active")
Set-TextRow $humidity 198 @("2026-01-28", "12:19:31", "12:00", "Bathroom", "87.8%", "Active")
Set-TextRow $humidity 199 @("2026-01-28", "12:19:32", "12:00", "Bathroom", "86.8%", "Active")
Set-TextRow $humidity 200 @("2026-01-28", "12:19:34", "12:00", "Bathroom", "87.7%", "Active")
Set-TextRow $humidity 201 @("2026-01-28", "12:19:36", "12:00", "Bathroom", "86.8%", "Active")
Set-TextRow $humidity 202 @("2026-01-28", "12:19:38", "12:00", "Bathroom", "86.2%", "Active")
Set-TextRow $humidity 203 @("2026-01-28", "12:19:40", "12:00", "Bathroom", "86.7%", "Active")
Set-TextRow $humidity 204 @("2026-01-28", "12:19:44", "12:00", "Bathroom", "87.6%", "Active")
Set-TextRow $humidity 205 @("2026-01-28", "12:19:52", "12:00", "Bathroom", "86.7%", "Active")
Set-TextRow $humidity 206 @("2026-01-28", "12:19:56", "12:00", "Bathroom", "87.6%", "Active")
Set-TextRow $humidity 207 @("2026-01-28", "12:20:00", "12:00", "Bathroom", "86.7%", "Active")
Set-TextRow $humidity 208 @("2026-01-28", "12:20:04", "12:00", "Bathroom", "87.6%", "Active")
Set-TextRow $humidity 209 @("2026-01-28", "12:20:08", "12:00", "Bathroom", "87.6%", "Active")
Set-TextRow $humidity 210 @("2026-01-28", "12:20:16", "12:00", "Bathroom", "87.6%", "Active")
$humidityRange.Style = "Normal"

$temperature = $wb.Worksheets.Item("Temperature")
$temperatureRange = $temperature.Range("A191:F210")
$temperatureRange.NumberFormat = "@"
Set-TextRow $temperature 191 @("2026-01-28", "12:19:19", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 192 @("2026-01-28", "12:19:22", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 193 @("2026-01-28", "12:19:23", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 194 @("2026-01-28", "12:19:25", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 195 @("2026-01-28", "12:19:27", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 196 @("2026-01-28", "12:19:28", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 197 @("2026-01-28", "12:19:30", "12:00", "Bathroom", "22.9C", "Active")
Set-TextRow $temperature 198 @("2026-01-28", "12:19:31", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 199 @("2026-01-28", "12:19:33", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 200 @("2026-01-28", "12:19:35", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 201 @("2026-01-28", "12:19:36", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 202 @("2026-01-28", "12:19:38", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 203 @("2026-01-28", "12:19:41", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 204 @("2026-01-28", "12:19:45", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 205 @("2026-01-28", "12:19:53", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 206 @("2026-01-28", "12:19:57", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 207 @("2026-01-28", "12:20:01", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 208 @("2026-01-28", "12:20:05", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 209 @("2026-01-28", "12:20:09", "12:00", "Bathroom", "23.0C", "Active")
Set-TextRow $temperature 210 @("2026-01-28", "12:20:17", "12:00", "Bathroom", "23.0C", "Active")
$temperatureRange.Style = "Normal"

